$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the "Meta description: ..." paragraph that currently sits
#    right after the "Play Black Knight Slot for Free - Read Our
#    Review" Heading1 paragraph (paragraph #2).
# ------------------------------------------------------------------
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# ------------------------------------------------------------------
# 2. Insert a new paragraph - bold "Play Black Knight Slot for Free -
#    Read Our Review" - right before the final ("Please create a
#    cartoon-style ...") paragraph.
# ------------------------------------------------------------------
$paraCount = $d.Paragraphs.Count
$secondToLast = $d.Paragraphs.Item($paraCount - 1)
$newParaRange = $secondToLast.Range.InsertParagraphAfter()

$insertedPara = $d.Paragraphs.Item($paraCount)
$xmlFragment = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Black Knight Slot for Free - Read Our Review</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertedPara.Range.InsertXML($xmlFragment)

# ------------------------------------------------------------------
# 3. Replace the text of the last paragraph (previously the AI image
#    prompt, still italic) with the meta-description text.
# ------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastRange = $d.Range($lastPara.Range.Start, $lastPara.Range.End)
$lastRange.Text = "Discover Black Knight slot machine by WMS Gaming. Check pros & cons, bet limits, and free spins feature. Play now for free at our site."
